$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-15 (columns A-F), as per the diff.
$data = @(
    @(45953, 7224.69439304514, 6626.39699882839, 9620, 3084.58590258218, 3.79095422544037),
    @(45954, 7466.88524738323, 6375.78272644713, 3620, 3510.01537007393, 261.074920688378),
    @(45955, 3251.19195569176, 4829.47250770238, 3620, 3786.78615303061, 208.177444197208),
    @(45956, 3096.60264544039, 4643.41151374032, 3620, 3729.41693478243, 198.034518688448),
    @(45957, 8643.4557607038, 7939.30652176617, 3620, 4438.77363760911, 364.920006640636),
    @(45958, 7224.82758554897, 7265.52152705163, 3620, 3510.48993179286, 298.167144118521),
    @(45959, 7224.82758554897, 7450.76546262333, 3620, 3510.48993179286, 305.885641434008),
    @(45960, 7224.82758554897, 7375.34445168465, 3620, 3510.48993179286, 302.743099311563),
    @(45961, 6107.35191671871, 5804.91029276995, 3620, 3445.87224974344, 234.615939271391),
    @(45962, 2557.16655997904, 4083.27584532022, 9004, 3794.11847801793, -46.941903194244),
    @(45963, 2557.16655997904, 4197.64852787918, 9004, 3794.11847801793, -42.1763747542872),
    @(45964, 8571.71877057272, 7238.09493097941, 9004, 4682.6673025921, 121.531759732146),
    @(45965, 8571.71877057272, 7115.29807806796, 9004, 4682.6673025921, 116.415224194169),
    @(45966, 8571.71877057272, 7276.47725853153, 9004, 4682.6673025921, 123.131023380151)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
}

$wb.Save()
